# gen gabors to test range model confidence
#
# Rename the original sheet to "Model 1", add a second sheet "Model 2" that
# documents the new "range model" (trained across a spread of tilts/contrasts
# with a confidence threshold), and update the training-data caption on
# "Model 1" to record the threshold that was used for that run too.

$wb = $excel.ActiveWorkbook

# --- Sheet1 -> "Model 1" -------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Model 1"

# --- add "Model 2" right after "Model 1" ---------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Model 2"

# ===========================================================================
# Populate "Model 2" (order matters: it determines shared-string order)
# ===========================================================================

# Results blurb for the range model (wrapped, 3-line caption) - row 5
$ws2.Range("C5").Value2 = "High Conf: 0.38`nTf Acc: 0.61`nMy Acc: 0.66"
$ws2.Range("C5").WrapText = $true
$ws2.Rows.Item(5).RowHeight = 45

# Small legend under the blurb
$ws2.Range("B7").Value2 = "tilts"
$ws2.Range("C7").Value2 = "contrasts 0.3, 0.45, 1"

# Header / title for the sheet
$ws2.Range("C2").Value2 = "Range model"
$ws2.Range("C2").Font.Bold = $true
$ws2.Range("C2").Font.Size = 11

# Update "Model 1"'s training caption to note the confidence threshold used
$ws1.Range("C3").Value2 = "Trained on: 20,000 images with tilt 2.26 + contrast 1, threshold = 2"

# Training-data caption for the new range model
$ws2.Range("C3").Value2 = "Trained on: 18,000 images with tilts [0.1, 0.2, 0.4, 0.8, 1.6, 3.2] & contrasts [0.3, 0.45, 1], threshold = 0.5"

# List of tilts used to generate the gabor test images
$ws2.Range("B8").Value2 = 0.1
$ws2.Range("B9").Value2 = 0.2
$ws2.Range("B10").Value2 = 0.4
$ws2.Range("B11").Value2 = 0.8
$ws2.Range("B12").Value2 = 1.6
$ws2.Range("B13").Value2 = 3.2

# Match column widths used elsewhere in the workbook
$ws2.Columns.Item(3).ColumnWidth = 13

# ===========================================================================
# Selections / active views
# ===========================================================================
$ws1.Range("C3").Select()

$ws2.Activate()
$ws2.Range("C9").Select()
